# Update average_county_temperature (I), worst_ashp_cop (N) and best_ashp_cop (O)
# values with refreshed NOAA temperature data for NAICS 312140 longform electrification
# options dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> average_county_temperature (column I)
$tempUpdates = @{
    2  = 13.75752314814816
    5  = 13.75752314814816
    6  = 3.38888888888889
    8  = 12.93898809523811
    9  = 12.93898809523811
    11 = 19.79629629629628
    12 = 19.79629629629628
    13 = 13.75752314814816
    14 = 3.38888888888889
    15 = 3.38888888888889
}

foreach ($row in $tempUpdates.Keys) {
    $ws.Range("I$row").Value = $tempUpdates[$row]
}

# Map of row -> worst_ashp_cop (column N) for rows that carry a COP value
$worstCopUpdates = @{
    2  = 1.722630989917367
    5  = 1.722630989917367
    6  = 1.578134831460674
    11 = 1.819666609086197
    13 = 1.722630989917367
    14 = 1.578134831460674
}

foreach ($row in $worstCopUpdates.Keys) {
    $ws.Range("N$row").Value = $worstCopUpdates[$row]
}

# Map of row -> best_ashp_cop (column O) for rows that carry a COP value
$bestCopUpdates = @{
    2  = 1.865269081797952
    5  = 1.865269081797952
    6  = 1.695036674816626
    11 = 1.981148790245761
    13 = 1.865269081797952
    14 = 1.695036674816626
}

foreach ($row in $bestCopUpdates.Keys) {
    $ws.Range("O$row").Value = $bestCopUpdates[$row]
}
